$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 16

$ws.Cells.Item($row, 1).Value = 42625.886678240742

$ws.Cells.Item($row, 2).Value = -16
$ws.Cells.Item($row, 3).Value = 55
$ws.Cells.Item($row, 4).Value = 44
$ws.Cells.Item($row, 5).Value = 33
$ws.Cells.Item($row, 6).Value = 66
$ws.Cells.Item($row, 7).Value = 7117
$ws.Cells.Item($row, 8).Value = 13196
$ws.Cells.Item($row, 9).Value = 1190
$ws.Cells.Item($row, 10).Value = 191
$ws.Cells.Item($row, 11).Value = 153
$ws.Cells.Item($row, 12).Value = 2
$ws.Cells.Item($row, 13).Value = 4
$ws.Cells.Item($row, 14).Value = "Bag"
